$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the explicit center-alignment style override on C17:C47 so they
# fall back to the column-level style (same visual result, but matches the
# target OOXML which drops the redundant s="3" attribute on these cells).
$ws.Range("C17:C47").Style = "Normal"

# Replace the old tail rows (48-52) and append new rows (53-62) with the
# newly uploaded product data. Column C (Price) uses right alignment.
$ws.Range("A48").Value = 'VEER/A-38510'
$ws.Range("B48").Value = 'NICKER SUIT'
$ws.Range("C48").Value = 400
$ws.Range("C48").HorizontalAlignment = -4152
$ws.Range("D48").Value = 'VEER/A-38510 - 22-26 - 3PCS BOX - PC/COTTON - 400'
$ws.Range("E48").Value = 'https://i.postimg.cc/DzSpdRmP/Whats-App-Image-2025-05-28-at-15-45-04-1.jpg'

$ws.Range("A49").Value = 'CARELINE/A-4153'
$ws.Range("B49").Value = 'NICKER SUIT'
$ws.Range("C49").Value = 325
$ws.Range("C49").HorizontalAlignment = -4152
$ws.Range("D49").Value = 'CARELINE/A-4153 - 22-26 - 3PCS BOX - PC/COTTON - 325'
$ws.Range("E49").Value = 'https://i.postimg.cc/Y9gXTg64/Whats-App-Image-2025-05-28-at-15-45-05.jpg'

$ws.Range("A50").Value = 'BENTLY/A-2175/4'
$ws.Range("B50").Value = 'NICKER SUIT'
$ws.Range("C50").Value = 319
$ws.Range("C50").HorizontalAlignment = -4152
$ws.Range("D50").Value = 'BENTLY/A-2175/4 - 22-26 - 3PCS BOX - PC/COTTON - 319'
$ws.Range("E50").Value = 'https://i.postimg.cc/prH0vTCs/Whats-App-Image-2025-05-28-at-15-45-05-1.jpg'

$ws.Range("A51").Value = 'KIDSTYLE/A-3066'
$ws.Range("B51").Value = 'TEXTILE SUIT'
$ws.Range("C51").Value = 345
$ws.Range("C51").HorizontalAlignment = -4152
$ws.Range("D51").Value = 'KIDSTYLE/A-3066 - 22-26 - 3PCS BOX - PC/COTTON - 345'
$ws.Range("E51").Value = 'https://i.postimg.cc/6p3MtFVx/Whats-App-Image-2025-05-28-at-15-45-06.jpg'

$ws.Range("A52").Value = 'BARRONBOY''S/A-7223'
$ws.Range("B52").Value = 'NICKER SUIT'
$ws.Range("C52").Value = 360
$ws.Range("C52").HorizontalAlignment = -4152
$ws.Range("D52").Value = 'BARRONBOY''S/A-7223 - 22-26 - 3PCS BOX - PC/COTTON - 360'
$ws.Range("E52").Value = 'https://i.postimg.cc/1Xrv3QYs/Whats-App-Image-2025-05-28-at-15-45-06-1.jpg'

$ws.Range("A53").Value = 'BARRONBOY''S/A-3306'
$ws.Range("B53").Value = 'CORD-SET'
$ws.Range("C53").Value = 180
$ws.Range("C53").HorizontalAlignment = -4152
$ws.Range("D53").Value = 'BARRONBOY''S/A-3306 - 18 - 3PCS BOX - PC/COTTON - 180'
$ws.Range("E53").Value = 'https://i.postimg.cc/5NfPkX60/Whats-App-Image-2025-05-28-at-15-45-06-2.jpg'

$ws.Range("A54").Value = 'CARELINE/A-4191'
$ws.Range("B54").Value = 'TEXTILE SUIT'
$ws.Range("C54").Value = 315
$ws.Range("C54").HorizontalAlignment = -4152
$ws.Range("D54").Value = 'CARELINE/A-4191 - 22-26 - 3PCS BOX - PC/COTTON - 315'
$ws.Range("E54").Value = 'https://i.postimg.cc/Gm5q2zz8/Whats-App-Image-2025-05-28-at-15-45-07.jpg'

$ws.Range("A55").Value = 'HEY DUDE(CAL)/A-1321'
$ws.Range("B55").Value = 'DANGRI SUIT'
$ws.Range("C55").Value = 288
$ws.Range("C55").HorizontalAlignment = -4152
$ws.Range("D55").Value = 'HEY DUDE(CAL)/A-1321 - 0.1.2 - 3PCS BOX - PC/COTTON - 288'
$ws.Range("E55").Value = 'https://i.postimg.cc/3J9LM3vZ/Whats-App-Image-2025-05-28-at-15-45-07-1.jpg'

$ws.Range("A56").Value = 'KIDSTYLE/A-3069'
$ws.Range("B56").Value = 'TEXTILE SUIT'
$ws.Range("C56").Value = 345
$ws.Range("C56").HorizontalAlignment = -4152
$ws.Range("D56").Value = 'KIDSTYLE/A-3069 - 22-26 - 3PCS BOX - PC/COTTON - 345'
$ws.Range("E56").Value = 'https://i.postimg.cc/pTNCQTwP/Whats-App-Image-2025-05-28-at-15-45-08.jpg'

$ws.Range("A57").Value = 'TOP GEAR/A-3322'
$ws.Range("B57").Value = 'NICKER SUIT'
$ws.Range("C57").Value = 369
$ws.Range("C57").HorizontalAlignment = -4152
$ws.Range("D57").Value = 'TOP GEAR/A-3322 - 22/26 - 3PCS BOX - PC/COTTON - 369'
$ws.Range("E57").Value = 'https://i.postimg.cc/yxsT5G65/Whats-App-Image-2025-05-28-at-15-45-09.jpg'

$ws.Range("A58").Value = 'HELLO KITTY/A-1380'
$ws.Range("B58").Value = 'DANGRI SUIT'
$ws.Range("C58").Value = 230
$ws.Range("C58").HorizontalAlignment = -4152
$ws.Range("D58").Value = 'HELLO KITTY/A-1380 - 0.1.2. - 3PCS BOX - PC/COTTON - 230'
$ws.Range("E58").Value = 'https://i.postimg.cc/FRjZTvv3/Whats-App-Image-2025-05-28-at-15-45-11.jpg'

$ws.Range("A59").Value = 'BARRONBOY''S/A-8680'
$ws.Range("B59").Value = 'TEXTILE SUIT'
$ws.Range("C59").Value = 396
$ws.Range("C59").HorizontalAlignment = -4152
$ws.Range("D59").Value = 'BARRONBOY''S/A-8680 - 28/32 - 3PCS BOX - PC/COTTON - 396'
$ws.Range("E59").Value = 'https://i.postimg.cc/Kjf5y5RY/Whats-App-Image-2025-05-28-at-15-45-12-1.jpg'

$ws.Range("A60").Value = 'B.M /A-3621'
$ws.Range("B60").Value = 'NICKER SUIT'
$ws.Range("C60").Value = 220
$ws.Range("C60").HorizontalAlignment = -4152
$ws.Range("D60").Value = 'B.M /A-3621 - 20 - 3PCS BOX - PC/COTTON - 220'
$ws.Range("E60").Value = 'https://i.postimg.cc/GhMj7zW7/Whats-App-Image-2025-05-28-at-15-45-13.jpg'

$ws.Range("A61").Value = 'VEER/A-83113'
$ws.Range("B61").Value = 'NICKER SUIT'
$ws.Range("C61").Value = 438
$ws.Range("C61").HorizontalAlignment = -4152
$ws.Range("D61").Value = 'VEER/A-83113 - 22-26 - 3PCS BOX - PC/COTTON - 438'
$ws.Range("E61").Value = 'https://i.postimg.cc/J75xJRPW/Whats-App-Image-2025-05-28-at-15-45-13-1.jpg'

$ws.Range("A62").Value = 'HEY DUDE(CAL)/A-1386'
$ws.Range("B62").Value = 'DANGRI SUIT'
$ws.Range("C62").Value = 246
$ws.Range("C62").HorizontalAlignment = -4152
$ws.Range("D62").Value = 'HEY DUDE(CAL)/A-1386 - 0.1.2 - 3PCS BOX - PC/COTTON - 246'
$ws.Range("E62").Value = 'https://i.postimg.cc/BnXBdHsq/Whats-App-Image-2025-05-28-at-15-45-14.jpg'

# Trailing blank row with the same right-aligned price style, matching the
# source workbook (row 63, only C has formatting, no values).
$ws.Range("C63").HorizontalAlignment = -4152

# Restore the selection/scroll state from the diff.
$excel.ActiveWindow.ScrollRow = 38
$ws.Range("D52").Select()
